$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2: replace the two-column (Subject/Attendance?) table with a
#     single "Roll No." column listing 1..60 ---
$ws2.Cells.Clear()
$ws2.Range("A1").Value = "Roll No."
for ($i = 1; $i -le 60; $i++) {
    $ws2.Cells.Item($i + 1, 1).Value = $i
}

# --- Selections / active tab: Sheet2 was tabSelected before, now Sheet1 is ---
$ws2.Range("G11").Select() | Out-Null
$ws1.Range("F11").Select() | Out-Null
$ws1.Activate() | Out-Null
